# Apply the "Ajout draft mapping" edit:
#  - Update the generation Date on the Metadata sheet
#  - Add a new "Mapping: Spécification métier vers l'extension ROR
#    NbTemporarySocialHelpPlace" column (AL) on the Elements sheet, with a
#    value only on the Extension.value[x] row (nbPlaceAideSocialTemporaire)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date value (B8, next to the "Date" label) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# --- Elements sheet: new mapping column AL ---
$wsElem = $wb.Worksheets.Item("Elements")

# Copy formatting from the existing last column (AK) so the new column (AL)
# uses the same header/data styles, then fill in the actual content.
$wsElem.Range("AK1").Copy()
$wsElem.Range("AL1").PasteSpecial(-4122)
$wsElem.Range("AK2:AK6").Copy()
$wsElem.Range("AL2:AL6").PasteSpecial(-4122)

# Header cell (row 1)
$wsElem.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR NbTemporarySocialHelpPlace"

# Data cells (rows 2-6) - only the Extension.value[x] row carries a mapping
$wsElem.Range("AL2").Value = " "
$wsElem.Range("AL3").Value = " "
$wsElem.Range("AL4").Value = " "
$wsElem.Range("AL5").Value = " "
$wsElem.Range("AL6").Value = "nbPlaceAideSocialTemporaire"

# Column width for the new column
$wsElem.Columns.Item(38).ColumnWidth = 84.69010416666667
